# Add a new test-file entry (BadanieControllerTests.cs) to the report,
# as a new row-26 entry in the "K/L/M" (4th person) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for row 26, columns K (date), L (file name), M (lines).
$ws.Range("K26").Value = 45793
$ws.Range("L26").Value = "BadanieControllerTests.cs"
$ws.Range("M26").Value = 55

# Match the date formatting used by the rest of the column (copy K25's
# format, e.g. date number format + style, onto K26) without introducing
# a brand-new style entry.
$ws.Range("K25").Copy()
$ws.Range("K26").PasteSpecial(-4122)

# Recalculate so the dependent SUM()/ratio formulas (M3, D4, G4, J4, M4, P4)
# pick up the new value.
$excel.Calculate()
